# Updates cryptos list data to reflect latest prices/volume changes
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = '41.057.21'
$ws.Cells.Item(2, 5).NumberFormat = "@"
$ws.Cells.Item(2, 5).Value = '  -1.21%  '
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = '2.171.18'
$ws.Cells.Item(3, 5).NumberFormat = "@"
$ws.Cells.Item(3, 5).Value = '  -2.13%  '
$ws.Cells.Item(4, 5).NumberFormat = "@"
$ws.Cells.Item(4, 5).Value = '  -0.10%  '
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '247.15'
$ws.Cells.Item(5, 5).NumberFormat = "@"
$ws.Cells.Item(5, 5).Value = '  -1.18%  '
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '0.614'
$ws.Cells.Item(6, 5).NumberFormat = "@"
$ws.Cells.Item(6, 5).Value = '  -1.80%  '
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '66.36'
$ws.Cells.Item(7, 5).NumberFormat = "@"
$ws.Cells.Item(7, 5).Value = '  -5.68%  '
$ws.Cells.Item(8, 5).NumberFormat = "@"
$ws.Cells.Item(8, 5).Value = '  -0.07%  '
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '0.564'
$ws.Cells.Item(9, 5).NumberFormat = "@"
$ws.Cells.Item(9, 5).Value = '  -1.19%  '
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '59.61'
$ws.Cells.Item(10, 5).NumberFormat = "@"
$ws.Cells.Item(10, 5).Value = '  +1.93%  '
$ws.Cells.Item(11, 5).NumberFormat = "@"
$ws.Cells.Item(11, 5).Value = '  -3.39%  '
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '35.43'
$ws.Cells.Item(12, 5).NumberFormat = "@"
$ws.Cells.Item(12, 5).Value = '  -14.38%  '
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '0.104'
$ws.Cells.Item(13, 5).NumberFormat = "@"
$ws.Cells.Item(13, 5).Value = '  -1.56%  '
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '2.488.70'
$ws.Cells.Item(15, 5).NumberFormat = "@"
$ws.Cells.Item(15, 5).Value = '  -2.26%  '
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '0.854'
$ws.Cells.Item(16, 5).NumberFormat = "@"
$ws.Cells.Item(16, 5).Value = '  +0.25%  '
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = '14.24'
$ws.Cells.Item(17, 5).NumberFormat = "@"
$ws.Cells.Item(17, 5).Value = '  -3.94%  '
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '2.178.05'
$ws.Cells.Item(18, 5).NumberFormat = "@"
$ws.Cells.Item(18, 5).Value = '  -1.86%  '
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '40.958.90'
$ws.Cells.Item(19, 5).NumberFormat = "@"
$ws.Cells.Item(19, 5).Value = '  -1.23%  '
$ws.Cells.Item(20, 5).NumberFormat = "@"
$ws.Cells.Item(20, 5).Value = '  -2.57%  '
$ws.Cells.Item(21, 2).NumberFormat = "@"
$ws.Cells.Item(21, 2).Value = 'Litecoin'
$ws.Cells.Item(21, 3).NumberFormat = "@"
$ws.Cells.Item(21, 3).Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '71.46'
$ws.Cells.Item(21, 5).NumberFormat = "@"
$ws.Cells.Item(21, 5).Value = '  -1.36%  '
$ws.Cells.Item(22, 2).NumberFormat = "@"
$ws.Cells.Item(22, 2).Value = 'Uniswap'
$ws.Cells.Item(22, 3).NumberFormat = "@"
$ws.Cells.Item(22, 3).Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '6.07'
$ws.Cells.Item(22, 5).NumberFormat = "@"
$ws.Cells.Item(22, 5).Value = '  -1.82%  '
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '229.45'
$ws.Cells.Item(23, 5).NumberFormat = "@"
$ws.Cells.Item(23, 5).Value = '  -1.81%  '
$ws.Cells.Item(24, 5).NumberFormat = "@"
$ws.Cells.Item(24, 5).Value = '  -6.76%  '
$ws.Cells.Item(25, 5).NumberFormat = "@"
$ws.Cells.Item(25, 5).Value = '  +0.02%  '
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '11.30'
$ws.Cells.Item(26, 5).NumberFormat = "@"
$ws.Cells.Item(26, 5).Value = '  +8.02%  '
$ws.Cells.Item(27, 5).NumberFormat = "@"
$ws.Cells.Item(27, 5).Value = '  -4.99%  '
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '2.43'
$ws.Cells.Item(28, 5).NumberFormat = "@"
$ws.Cells.Item(28, 5).Value = '  -2.83%  '
$ws.Cells.Item(29, 5).NumberFormat = "@"
$ws.Cells.Item(29, 5).Value = '  -5.45%  '
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '168.54'
$ws.Cells.Item(30, 5).NumberFormat = "@"
$ws.Cells.Item(30, 5).Value = '  -1.55%  '
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '2.00'
$ws.Cells.Item(31, 5).NumberFormat = "@"
$ws.Cells.Item(31, 5).Value = '  -8.58%  '
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '20.17'
$ws.Cells.Item(32, 5).NumberFormat = "@"
$ws.Cells.Item(32, 5).Value = '  -1.76%  '
$ws.Cells.Item(33, 5).NumberFormat = "@"
$ws.Cells.Item(33, 5).Value = '  +2.03%  '
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '5.64'
$ws.Cells.Item(34, 5).NumberFormat = "@"
$ws.Cells.Item(34, 5).Value = '  +2.61%  '
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '0.0746'
$ws.Cells.Item(35, 5).NumberFormat = "@"
$ws.Cells.Item(35, 5).Value = '  +4.19%  '
$ws.Cells.Item(36, 5).NumberFormat = "@"
$ws.Cells.Item(36, 5).Value = '  -2.42%  '
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '4.54'
$ws.Cells.Item(37, 5).NumberFormat = "@"
$ws.Cells.Item(37, 5).Value = '  -2.54%  '
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '3.99'
$ws.Cells.Item(38, 5).NumberFormat = "@"
$ws.Cells.Item(38, 5).Value = '  +1.34%  '
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '24.25'
$ws.Cells.Item(39, 5).NumberFormat = "@"
$ws.Cells.Item(39, 5).Value = '  -6.32%  '
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '0.0301'
$ws.Cells.Item(40, 5).NumberFormat = "@"
$ws.Cells.Item(40, 5).Value = '  +3.32%  '
$ws.Cells.Item(41, 5).NumberFormat = "@"
$ws.Cells.Item(41, 5).Value = '  -4.10%  '
$ws.Cells.Item(42, 5).NumberFormat = "@"
$ws.Cells.Item(42, 5).Value = '  -7.51%  '
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '4.85'
$ws.Cells.Item(43, 5).NumberFormat = "@"
$ws.Cells.Item(43, 5).Value = '  -0.99%  '
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '60.36'
$ws.Cells.Item(44, 5).NumberFormat = "@"
$ws.Cells.Item(44, 5).Value = '  -12.11%  '
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '11.17'
$ws.Cells.Item(45, 5).NumberFormat = "@"
$ws.Cells.Item(45, 5).Value = '  -5.22%  '
$ws.Cells.Item(46, 5).NumberFormat = "@"
$ws.Cells.Item(46, 5).Value = '  -7.76%  '
$ws.Cells.Item(47, 5).NumberFormat = "@"
$ws.Cells.Item(47, 5).Value = '  -2.85%  '
$ws.Cells.Item(48, 5).NumberFormat = "@"
$ws.Cells.Item(48, 5).Value = '  -0.10%  '
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '0.0987'
$ws.Cells.Item(49, 5).NumberFormat = "@"
$ws.Cells.Item(49, 5).Value = '  -2.30%  '
$ws.Cells.Item(50, 5).NumberFormat = "@"
$ws.Cells.Item(50, 5).Value = '  -0.50%  '
$ws.Cells.Item(51, 5).NumberFormat = "@"
$ws.Cells.Item(51, 5).Value = '  -3.08%  '

Write-Output "Update complete"
